$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: update per-week data with restock suggestion fields ---
$ws1.Range("B2").Value = "'2025-02-02"
$ws1.Range("L2").Value = 60.94
$ws1.Range("P2").Value = 1.03
$ws1.Range("Q2").Value = "Decline"

$ws1.Range("B3").Value = "'2025-02-09"
$ws1.Range("L3").Value = 59.94
$ws1.Range("P3").Value = 1
$ws1.Range("Q3").Value = "Decline"

$ws1.Range("B4").Value = "'2025-02-16"
$ws1.Range("L4").Value = 58.94
$ws1.Range("P4").Value = 0.82
$ws1.Range("Q4").Value = "Decline"

$ws1.Range("B5").Value = "'2025-02-23"
$ws1.Range("L5").Value = 57.94
$ws1.Range("P5").Value = 0.85
$ws1.Range("Q5").Value = "Decline"

$ws1.Range("B6").Value = "'2025-03-02"
$ws1.Range("L6").Value = 47.32
$ws1.Range("P6").Value = 0.82
$ws1.Range("Q6").Value = "Decline"

$ws1.Range("B7").Value = "'2025-03-09"
$ws1.Range("L7").Value = 46.32
$ws1.Range("P7").Value = 0.93
$ws1.Range("Q7").Value = "Decline"

$ws1.Range("B8").Value = "'2025-03-16"
$ws1.Range("L8").Value = 45.32
$ws1.Range("P8").Value = 0.9
$ws1.Range("Q8").Value = "Decline"

$ws1.Range("B9").Value = "'2025-03-23"
$ws1.Range("L9").Value = 44.32
$ws1.Range("P9").Value = 1.17
$ws1.Range("Q9").Value = "Decline"

$ws1.Range("B10").Value = "'2025-03-30"
$ws1.Range("L10").Value = 43.32
$ws1.Range("P10").Value = 1.13
$ws1.Range("Q10").Value = "Decline"

$ws1.Range("B11").Value = "'2025-04-06"
$ws1.Range("L11").Value = 42.32
$ws1.Range("P11").Value = 0.93
$ws1.Range("Q11").Value = "Decline"

$ws1.Range("B12").Value = "'2025-04-13"
$ws1.Range("L12").Value = 41.32
$ws1.Range("P12").Value = 1.1
$ws1.Range("Q12").Value = "Decline"

$ws1.Range("B13").Value = "'2025-04-20"
$ws1.Range("L13").Value = 41.4
$ws1.Range("P13").Value = 0.9399999999999999
$ws1.Range("Q13").Value = "Decline"

$ws1.Range("B14").Value = "'2025-04-27"
$ws1.Range("L14").Value = 40.4
$ws1.Range("P14").Value = 0.93
$ws1.Range("Q14").Value = "Decline"

$ws1.Range("B15").Value = "'2025-05-04"
$ws1.Range("L15").Value = 39.4
$ws1.Range("P15").Value = 1.07
$ws1.Range("Q15").Value = "Decline"

$ws1.Range("B16").Value = "'2025-05-11"
$ws1.Range("L16").Value = 38.4
$ws1.Range("P16").Value = 1.19
$ws1.Range("Q16").Value = "Decline"

$ws1.Range("B17").Value = "'2025-05-18"
$ws1.Range("L17").Value = 37.4
$ws1.Range("P17").Value = 1.08
$ws1.Range("Q17").Value = "Decline"

# Header: Q1 becomes "Lifecycle Stage" (replacing "Sales Volume Rank"),
# and the old R column ("Lifecycle Stage") is removed entirely.
$ws1.Range("Q1").Value = "Lifecycle Stage"
$ws1.Range("R1:R17").EntireColumn.Delete()

# --- Summary sheet: Max/Min Forecast Week are no longer available ---
$ws2.Range("B13").Value = "N/A"
$ws2.Range("B15").Value = "N/A"

